# Applies the diff:
#  - Row 10 ("P. point" / "Long point (up to 10 mtr.)") is removed; all rows
#    below it shift up by one (old row 19 is gone, dimension A1:I19 -> A1:I18).
#  - A handful of cells in the shifted rows (8, 9, 10, 11, 12, 13, 14, 16, 18)
#    get new values (quantities/labels/amounts recomputed for the new row set).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the whole row 10 - everything below shifts up one row.
$ws.Rows(10).Delete()

# 2) Row 8: "Qty executed upto date" revised down.
$ws.Range("C8").Value = 45

# 3) Row 9 (was "Short point (up to 3 mtr.)") becomes "Medium point (up to 6 mtr.)".
$ws.Range("C9").Value = 85
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3"
$ws.Range("E9").Value = "Medium point (up to 6 mtr.)"
$ws.Range("F9").Value = 472
$ws.Range("G9").NumberFormat = "@"
$ws.Range("G9").Value = "40120.00"

# 4) Row 10 (was row 11, the long "Rewiring of 3/5 pin..." text row): qty tweak only.
$ws.Range("C10").Value = 12

# 5) Row 11 (was row 12, "On board"): qty + amount tweak.
$ws.Range("C11").Value = 88
$ws.Range("G11").NumberFormat = "@"
$ws.Range("G11").Value = "11968.00"

# 6) Row 12 (was row 13, "P & F ISI marked..." switch row) becomes the "Total" row.
$ws.Range("A12").Value = ""
$ws.Range("C12").Value = 58
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8"
$ws.Range("E12").Value = "Total"
$ws.Range("F12").Value = 0
$ws.Range("G12").NumberFormat = "@"
$ws.Range("G12").Value = "0.00"

# 7) Row 13 ("Add Tender Premium"): qty tweak only.
$ws.Range("C13").Value = 5

# 8) Row 14 ("Grand Total"): qty tweak only.
$ws.Range("C14").Value = 26

# 9) Row 16 ("Grand Total Rs." summary row): new totals.
$ws.Range("G16").NumberFormat = "@"
$ws.Range("G16").Value = "52088.00"
$ws.Range("H16").NumberFormat = "@"
$ws.Range("H16").Value = "52088.00"

# 10) Row 18 ("NET PAYABLE AMOUNT Rs." row): new totals.
$ws.Range("G18").NumberFormat = "@"
$ws.Range("G18").Value = "52088.00"
$ws.Range("H18").NumberFormat = "@"
$ws.Range("H18").Value = "52088.00"
